$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E18").Value = "test3"
